# QA-TestCase-Auto-X-Alpha-004(): add 5 assertion points.
#
# The new assertions rely on a new "DealType" column (with value "FX Spot")
# inserted right after "TestCaseID" on the XAlphaDeals sheet, before the
# existing "Direction" column. After making this edit, the user left the
# XAlphaDeals sheet selected/active (instead of NitroXBots).

$wb = $excel.ActiveWorkbook

$wsDeals = $wb.Worksheets.Item("XAlphaDeals")

# Shift existing header/value cells B:G one column to the right (C:H) to make
# room for the new "DealType" column, and drop the old duplicate H column
# ("ProcessingStatus" / "processed") since I:N already carries that data.
$wsDeals.Range("H1").Value = $wsDeals.Range("G1").Value2
$wsDeals.Range("G1").Value = $wsDeals.Range("F1").Value2
$wsDeals.Range("F1").Value = $wsDeals.Range("E1").Value2
$wsDeals.Range("E1").Value = $wsDeals.Range("D1").Value2
$wsDeals.Range("D1").Value = $wsDeals.Range("C1").Value2
$wsDeals.Range("C1").Value = $wsDeals.Range("B1").Value2

$wsDeals.Range("H2").Value = $wsDeals.Range("G2").Value2
$wsDeals.Range("G2").Value = $wsDeals.Range("F2").Value2
$wsDeals.Range("F2").Value = $wsDeals.Range("E2").Value2
$wsDeals.Range("E2").Value = $wsDeals.Range("D2").Value2
$wsDeals.Range("D2").Value = $wsDeals.Range("C2").Value2
$wsDeals.Range("C2").Value = $wsDeals.Range("B2").Value2

# New DealType column + value (the 5 new assertion points).
$wsDeals.Range("B1").Value = "DealType"
$wsDeals.Range("B2").Value = "FX Spot"

# Selection / active-sheet bookkeeping matching the author's final state:
# XAlphaDeals (selecting the used range A1:XFD2) becomes the active sheet,
# replacing NitroXBots as the previously-selected tab.
$wsDeals.Range("A1:XFD2").Select()
